$wb = $excel.ActiveWorkbook

# The "Denmark" sheet is the last existing sheet and acts as the template
# for the three new market sheets (Russia, Finland, Hungary) that get
# appended right after it, in that order.
$denmark = $wb.Worksheets.Item("Denmark")

# --- Russia -----------------------------------------------------------
$denmark.Copy([System.Reflection.Missing]::Value, $denmark)
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B4").Value = "NGC-2929/T2903"
$russia.Range("B2").Value = "Russia Market"
[void]$russia.Range("A1:D16").Select()

# --- Finland ------------------------------------------------------------
$denmark.Copy([System.Reflection.Missing]::Value, $russia)
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2889"
$finland.Range("B2").Value = "Finland Market"
# Finland's printer list does not include the "MZX Communicator" row that
# the Denmark template has at row 12 - remove it so the row count matches.
$finland.Rows(12).Delete()
[void]$finland.Range("A1:D16").Select()

# --- Hungary ------------------------------------------------------------
$denmark.Copy([System.Reflection.Missing]::Value, $finland)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T2981"
$hungary.Range("B2").Value = "Hungary Market"

# Hungary ends up being the active/selected sheet.
[void]$hungary.Activate()
[void]$hungary.Range("H12").Select()
